# Updated test data for German,Czech market. Added test data for Belgium market
$wb = $excel.ActiveWorkbook

# Belgium: remove the stray "FAT-S" row (row 11) from the test data.
$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Activate()
$belgium.Rows.Item(11).Delete()
$belgium.Range("A11:XFD11").Select()

# Czech: remove the stray "FAT-S" row (row 11) from the test data.
$czech = $wb.Worksheets.Item("Czech")
$czech.Activate()
$czech.Rows.Item(11).Delete()
$czech.Range("A11:XFD11").Select()

# Focus returns to the Germany sheet.
$germany = $wb.Worksheets.Item("Germany")
$germany.Activate()
